$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H17 value
$ws.Range("H17").Value = 5

# Update D24 / E24 values and remove the green fill (switch to "no fill" style like G24/H24)
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 5
$ws.Range("D24:E24").Interior.Pattern = -4142  ## xlNone
$ws.Range("D24:E24").Interior.ColorIndex = -4142  ## xlNone

# Update G24 value
$ws.Range("G24").Value = 5

# Update the frozen-pane view / selection to match the scrolled-down state
$ws.Range("I17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.FreezePanes = $true
